# Manejo de errores consulta firestore
# Expands the stock table with additional tickers. Row 3 (previously NVDA.BA)
# is overwritten with AMZN.BA, new rows 4-12 are appended for the remaining
# tickers, and NVDA.BA's original data is preserved further down at row 9.
# PFE.BA (row 11) could not retrieve live price data (Firestore query
# failure), so its price/valuation columns (E-I) are left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing formatted cell (A3) as the template for the new
# bold/centered/bordered "index" cells in column A of the new rows.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 3: AMZN.BA
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "AMZN.BA"
$ws.Range("C3").Value = 1925
$ws.Range("D3").Value = 35
$ws.Range("E3").Value = 1930
$ws.Range("F3").Value = 67550
$ws.Range("G3").Value = 175
$ws.Range("H3").Value = 45.64
$ws.Range("I3").Value = 0.1199999999999974

# Row 4: COME.BA
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "COME.BA"
$ws.Range("C4").Value = 96.8
$ws.Range("D4").Value = 741
$ws.Range("E4").Value = 229.25
$ws.Range("F4").Value = 169874.25
$ws.Range("G4").Value = 98145.45
$ws.Range("H4").Value = 114.78
$ws.Range("I4").Value = 66.31

# Row 5: GGAL.BA
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "GGAL.BA"
$ws.Range("C5").Value = 3035
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 4110
$ws.Range("F5").Value = 28770
$ws.Range("G5").Value = 7525
$ws.Range("H5").Value = 19.44
$ws.Range("I5").Value = 5.090000000000002

# Row 6: GOOGL.BA
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "GOOGL.BA"
$ws.Range("C6").Value = 3998.94
$ws.Range("D6").Value = 17
$ws.Range("E6").Value = 4560
$ws.Range("F6").Value = 77520
$ws.Range("G6").Value = 9538.02
$ws.Range("H6").Value = 52.38
$ws.Range("I6").Value = 6.450000000000003

# Row 7: MELI.BA
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "MELI.BA"
$ws.Range("C7").Value = 18890.38
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = 20950
$ws.Range("F7").Value = 167600
$ws.Range("G7").Value = 16476.96
$ws.Range("H7").Value = 113.24
$ws.Range("I7").Value = 11.13

# Row 8: META.BA
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "META.BA"
$ws.Range("C8").Value = 10900
$ws.Range("D8").Value = 9
$ws.Range("E8").Value = 29650
$ws.Range("F8").Value = 266850
$ws.Range("G8").Value = 168750
$ws.Range("H8").Value = 180.3
$ws.Range("I8").Value = 114.02

# Row 9: NVDA.BA (moved here from the old row 3)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "NVDA.BA"
$ws.Range("C9").Value = 7163.88
$ws.Range("D9").Value = 21
$ws.Range("E9").Value = 7690
$ws.Range("F9").Value = 161490
$ws.Range("G9").Value = 11048.52
$ws.Range("H9").Value = 109.11
$ws.Range("I9").Value = 7.459999999999994

# Row 10: PAMP.BA
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "PAMP.BA"
$ws.Range("C10").Value = 2528
$ws.Range("D10").Value = 14
$ws.Range("E10").Value = 2690
$ws.Range("F10").Value = 37660
$ws.Range("G10").Value = 2268
$ws.Range("H10").Value = 25.45
$ws.Range("I10").Value = 1.539999999999999

# Row 11: PFE.BA - price lookup failed, leave E:I blank
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "PFE.BA"
$ws.Range("C11").Value = 9267.58
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""

# Row 12: VIST.BA
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "VIST.BA"
$ws.Range("C12").Value = 14377.5
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 22900
$ws.Range("F12").Value = 91600
$ws.Range("G12").Value = 34090
$ws.Range("H12").Value = 61.89
$ws.Range("I12").Value = 23.03
